# Re-sort the SAPTO rates table.
#
# The sheet lists SAPTO/SATO/PTO thresholds by financial year
# (fy_year, family_status_index). It had been sorted only by
# fy_year descending (and only over A2:K31); this updates the sort
# to cover the full data range A2:K49 and to additionally order by
# the sapto / sato flag columns (so rows with a source citation and
# rows belonging to the same scheme stay grouped together), matching
# Excel's new sortState:
#   descending by fy_year (A), then ascending by sapto (I), then
#   ascending by sato (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("A1:K49")
$keyYear   = $ws.Range("A2:A49")
$keySapto  = $ws.Range("I2:I49")
$keySato   = $ws.Range("J2:J49")

# Sort(Key1, Order1, Key2, Type2, Order2, Key3, Order3, Header)
$dataRange.Sort($keyYear, 2, $keySapto, $null, 1, $keySato, 1, 1)
